$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.67568533333333
$ws.Range("H2").Value = 83.027056
$ws.Range("I2").Value = 0.151580065893459
$ws.Range("J2").Value = 0.151580065893459
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 174.5186096564587
$ws.Range("R2").Value = 1570.667486908128
$ws.Range("S2").Value = 0.002055653839579177
$ws.Range("T2").Value = 0.002055653839579178
$ws.Range("G3").Value = 27.67568533333333
$ws.Range("H3").Value = 83.027056
$ws.Range("I3").Value = 0.151580065893459
$ws.Range("J3").Value = 0.151580065893459
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 5051.689403621293
$ws.Range("R3").Value = 45465.20463259164
$ws.Range("S3").Value = 0.05950382448815962
$ws.Range("T3").Value = 0.05950382448815962
$ws.Range("G4").Value = 27.67568533333333
$ws.Range("H4").Value = 83.027056
$ws.Range("I4").Value = 0.151580065893459
$ws.Range("J4").Value = 0.151580065893459
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 3525.782485269669
$ws.Range("R4").Value = 31732.04236742703
$ws.Range("S4").Value = 0.04153017444748695
$ws.Range("T4").Value = 0.04153017444748697
$ws.Range("G5").Value = 27.67568533333333
$ws.Range("H5").Value = 83.027056
$ws.Range("I5").Value = 0.151580065893459
$ws.Range("J5").Value = 0.151580065893459
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 537.5725488155805
$ws.Range("R5").Value = 4838.152939340224
$ws.Range("S5").Value = 0.006332064392447536
$ws.Range("T5").Value = 0.006332064392447536
$ws.Range("G6").Value = 27.67568533333333
$ws.Range("H6").Value = 83.027056
$ws.Range("I6").Value = 0.151580065893459
$ws.Range("J6").Value = 0.151580065893459
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 3579.112525357103
$ws.Range("R6").Value = 32212.01272821392
$ws.Range("S6").Value = 0.04215834872578569
$ws.Range("T6").Value = 0.04215834872578569
$ws.Range("I7").Value = 0.2439851776203359
$ws.Range("J7").Value = 0.243985177620336
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 280.9073457258773
$ws.Range("R7").Value = 2528.166111532896
$ws.Range("S7").Value = 0.003308806235301249
$ws.Range("T7").Value = 0.003308806235301249
$ws.Range("I8").Value = 0.2439851776203359
$ws.Range("J8").Value = 0.243985177620336
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.09577810315135514
$ws.Range("T8").Value = 0.09577810315135514
$ws.Range("I9").Value = 0.2439851776203359
$ws.Range("J9").Value = 0.243985177620336
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 5675.143765433019
$ws.Range("R9").Value = 51076.29388889717
$ws.Range("S9").Value = 0.06684749033092283
$ws.Range("T9").Value = 0.06684749033092283
$ws.Range("I10").Value = 0.2439851776203359
$ws.Range("J10").Value = 0.243985177620336
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 865.2835254655075
$ws.Range("R10").Value = 7787.551729189569
$ws.Range("S10").Value = 0.01019217036480642
$ws.Range("T10").Value = 0.01019217036480642
$ws.Range("I11").Value = 0.2439851776203359
$ws.Range("J11").Value = 0.243985177620336
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 5760.984467681938
$ws.Range("R11").Value = 51848.86020913744
$ws.Range("S11").Value = 0.06785860753795035
$ws.Range("T11").Value = 0.06785860753795034
$ws.Range("G12").Value = 54.059897
$ws.Range("H12").Value = 162.179691
$ws.Range("I12").Value = 0.2960867147735651
$ws.Range("J12").Value = 0.2960867147735651
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 340.893385257862
$ws.Range("R12").Value = 3068.040467320758
$ws.Range("S12").Value = 0.004015381498121704
$ws.Range("T12").Value = 0.004015381498121703
$ws.Range("G13").Value = 54.059897
$ws.Range("H13").Value = 162.179691
$ws.Range("I13").Value = 0.2960867147735651
$ws.Range("J13").Value = 0.2960867147735651
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 9867.64394617672
$ws.Range("R13").Value = 88808.79551559048
$ws.Range("S13").Value = 0.1162309292143029
$ws.Range("T13").Value = 0.1162309292143029
$ws.Range("G14").Value = 54.059897
$ws.Range("H14").Value = 162.179691
$ws.Range("I14").Value = 0.2960867147735651
$ws.Range("J14").Value = 0.2960867147735651
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 6887.035883751521
$ws.Range("R14").Value = 61983.32295376369
$ws.Range("S14").Value = 0.08112236159583366
$ws.Range("T14").Value = 0.08112236159583366
$ws.Range("G15").Value = 54.059897
$ws.Range("H15").Value = 162.179691
$ws.Range("I15").Value = 0.2960867147735651
$ws.Range("J15").Value = 0.2960867147735651
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 1050.059511407863
$ws.Range("R15").Value = 9450.535602670763
$ws.Range("S15").Value = 0.01236864578890096
$ws.Range("T15").Value = 0.01236864578890096
$ws.Range("G16").Value = 54.059897
$ws.Range("H16").Value = 162.179691
$ws.Range("I16").Value = 0.2960867147735651
$ws.Range("J16").Value = 0.2960867147735651
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 6991.207341094264
$ws.Range("R16").Value = 62920.86606984837
$ws.Range("S16").Value = 0.08234939667640591
$ws.Range("T16").Value = 0.08234939667640589
$ws.Range("G17").Value = 11.41370466666667
$ws.Range("H17").Value = 34.241114
$ws.Range("I17").Value = 0.0625129995743248
$ws.Range("J17").Value = 0.0625129995743248
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 71.97306391748134
$ws.Range("R17").Value = 647.7575752573321
$ws.Range("S17").Value = 0.0008477703637422521
$ws.Range("T17").Value = 0.000847770363742252
$ws.Range("G18").Value = 11.41370466666667
$ws.Range("H18").Value = 34.241114
$ws.Range("I18").Value = 0.0625129995743248
$ws.Range("J18").Value = 0.0625129995743248
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 2083.362714462485
$ws.Range("R18").Value = 18750.26443016236
$ws.Range("S18").Value = 0.02453991910462374
$ws.Range("T18").Value = 0.02453991910462374
$ws.Range("G19").Value = 11.41370466666667
$ws.Range("H19").Value = 34.241114
$ws.Range("I19").Value = 0.0625129995743248
$ws.Range("J19").Value = 0.0625129995743248
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 1454.064805300601
$ws.Range("R19").Value = 13086.58324770541
$ws.Range("S19").Value = 0.01712742214653845
$ws.Range("T19").Value = 0.01712742214653845
$ws.Range("G20").Value = 11.41370466666667
$ws.Range("H20").Value = 34.241114
$ws.Range("I20").Value = 0.0625129995743248
$ws.Range("J20").Value = 0.0625129995743248
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 221.6998146636062
$ws.Range("R20").Value = 1995.298331972456
$ws.Range("S20").Value = 0.002611401019893285
$ws.Range("T20").Value = 0.002611401019893284
$ws.Range("G21").Value = 11.41370466666667
$ws.Range("H21").Value = 34.241114
$ws.Range("I21").Value = 0.0625129995743248
$ws.Range("J21").Value = 0.0625129995743248
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 1476.058599495331
$ws.Range("R21").Value = 13284.52739545798
$ws.Range("S21").Value = 0.01738648693952707
$ws.Range("T21").Value = 0.01738648693952707
$ws.Range("G22").Value = 44.88488133333333
$ws.Range("H22").Value = 134.654644
$ws.Range("I22").Value = 0.2458350421383152
$ws.Range("J22").Value = 0.2458350421383153
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 283.0371494162746
$ws.Range("R22").Value = 2547.334344746472
$ws.Range("S22").Value = 0.003333893182431607
$ws.Range("T22").Value = 0.003333893182431607
$ws.Range("G23").Value = 44.88488133333333
$ws.Range("H23").Value = 134.654644
$ws.Range("I23").Value = 0.2458350421383152
$ws.Range("J23").Value = 0.2458350421383153
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 8192.912901105363
$ws.Range("R23").Value = 73736.21610994826
$ws.Range("S23").Value = 0.09650428052141961
$ws.Range("T23").Value = 0.09650428052141961
$ws.Range("G24").Value = 44.88488133333333
$ws.Range("H24").Value = 134.654644
$ws.Range("I24").Value = 0.2458350421383152
$ws.Range("J24").Value = 0.2458350421383153
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 5718.171982099698
$ws.Range("R24").Value = 51463.54783889728
$ws.Range("S24").Value = 0.06735431948212463
$ws.Range("T24").Value = 0.06735431948212464
$ws.Range("G25").Value = 44.88488133333333
$ws.Range("H25").Value = 134.654644
$ws.Range("I25").Value = 0.2458350421383152
$ws.Range("J25").Value = 0.2458350421383153
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 871.8439948651751
$ws.Range("R25").Value = 7846.595953786576
$ws.Range("S25").Value = 0.01026944610140217
$ws.Range("T25").Value = 0.01026944610140217
$ws.Range("G26").Value = 44.88488133333333
$ws.Range("H26").Value = 134.654644
$ws.Range("I26").Value = 0.2458350421383152
$ws.Range("J26").Value = 0.2458350421383153
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 5804.663517611676
$ws.Range("R26").Value = 52241.97165850508
$ws.Range("S26").Value = 0.06837310285093724
$ws.Range("T26").Value = 0.06837310285093723
